$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (column E) entirely, shifting the
# remaining columns (reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) one place to the
# left.
$ws.Range("E:E").Delete()
